$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 822 (shifts existing rows 822:863 down to 823:864)
$ws.Rows.Item(822).Insert()

# Populate the newly inserted row 822 with its data.
# Column A holds a date-like string that must stay plain text (not be
# auto-converted to a date serial number by the COM Value setter), so
# force a text number format before assignment, then strip the format
# override back off so the cell matches the unstyled data cells around it.
$ws.Cells.Item(822, 1).NumberFormat = "@"
$ws.Cells.Item(822, 1).Value = "2026/02/20"
$ws.Cells.Item(822, 1).ClearFormats()

$ws.Cells.Item(822, 2).Value = "金"
$ws.Cells.Item(822, 3).Value = 20
$ws.Cells.Item(822, 4).Value = 67
